$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing Call Date / Due Date values (rows 2 and 3) ---
# Row 2: Call Date 2022-12-10 -> 2022-01-10 ; Due Date 2022-12-19 -> 2022-01-19
$ws.Range("E2").Value = 44571
$ws.Range("F2").Value = 44580

# Row 3: Call Date 2022-12-10 -> 2022-06-10 ; Due Date 2022-12-19 -> 2022-06-19
$ws.Range("E3").Value = 44722
$ws.Range("F3").Value = 44731

# Row 4 (Call Date / Due Date) stay as-is (2022-12-10 / 2022-12-19)

# --- New header columns: From Currency, To Currency, Exchange Rate, As Of ---
$ws.Range("J1").Value = "From Currency"
$ws.Range("K1").Value = "To Currency"
$ws.Range("L1").Value = "Exchange Rate "
$ws.Range("M1").Value = "As Of"

# --- Row 2 FX data ---
$ws.Range("J2").Value = "USD"
$ws.Range("K2").Value = "INR"
$ws.Range("L2").Value = 80
$ws.Range("M2").Value = 44571

# --- Row 3 FX data ---
$ws.Range("J3").Value = "USD"
$ws.Range("K3").Value = "INR"
$ws.Range("L3").Value = 81
$ws.Range("M3").Value = 44722

# --- Row 4 FX data ---
$ws.Range("J4").Value = "USD"
$ws.Range("K4").Value = "INR"
$ws.Range("L4").Value = 82
$ws.Range("M4").Value = 44905

# Apply the same date style (dd/mm/yy, like column E) to the new "As Of" column
$ws.Range("M2:M4").NumberFormat = $ws.Range("E2").NumberFormat

# New column M width (best-fit like Excel autosizing)
$ws.Columns.Item(13).ColumnWidth = 9.25

# Update the active selection to A5 (was B5)
$ws.Range("A5").Select() | Out-Null
